# Helper: VBA-style RGB() -> long color value (R + G*256 + B*65536)
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1. Refresh the stale cached "datetimeFigureOut" footer field text
#    (8/9/2021 -> 9/3/2021) on every slide layout, the slide master,
#    and the notes master. msoPlaceholderDate = 16.
# ---------------------------------------------------------------
$newDate = "9/3/2021"

function Update-DatePlaceholder($shapes) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if ($shp.Type -eq 14) {
      if ($shp.PlaceholderFormat.Type -eq 16) {
        $shp.TextFrame.TextRange.Text = $newDate
      }
    }
  }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
  $layout = $master.CustomLayouts.Item($L)
  Update-DatePlaceholder $layout.Shapes
}

# NOTE: the notes master's own "Date Placeholder 2" (field id
# {D34CEA9E-...}) should also move from 8/9/2021 -> 9/3/2021, but this
# host's COM shim mis-seats NotesMaster shape writes onto the slide
# master (clobbering an unrelated placeholder), so it is intentionally
# left untouched here rather than risk corrupting the slide master.

# ---------------------------------------------------------------
# 2. Capitalize the "a" / "b" axis labels on the scatter-plot slide
# ---------------------------------------------------------------
$slide = $p.Slides.Item(1)
$slide.Shapes.Item("TextBox 51").TextFrame.TextRange.Text = "A"
$slide.Shapes.Item("TextBox 24").TextFrame.TextRange.Text = "B"

# ---------------------------------------------------------------
# 3. Lighten the "Top 100" / "Bottom 100" legend swatch colors
# ---------------------------------------------------------------
$slide.Shapes.Item("Rectangle 60").Fill.ForeColor.RGB = RGB 0xE8 0xF3 0xFF
$slide.Shapes.Item("Rectangle 71").Fill.ForeColor.RGB = RGB 0xFF 0xF5 0xE5
